# feat: adição de funcionalidades e algumas correções
#
# Extends Sheet1 from A1:K5 to A1:AJ9:
#  - adds 24 new header columns (L1:AJ1) with the same header style as A1:K1
#  - extends the existing data rows (2-5) with blank cells in the new columns
#  - appends 4 new data rows (6-9) with a mix of numeric / text values and
#    blank cells across the whole A:AJ span

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header row: copy the formatting of A1 (bold, border, centered) onto
#        the new header cells L1:AJ1, then fill in their text values. ------
$ws.Range("A1").Copy()
$ws.Range("L1:AJ1").PasteSpecial(-4122)

$ws.Range("L1").Value = "Nº Proposta:"
$ws.Range("M1").Value = "Unnamed: 12"
$ws.Range("N1").Value = "Nº Pedido:"
$ws.Range("O1").Value = "Nº Pedido de Vendas:"
$ws.Range("P1").Value = "Data Entrada do Pedido:"
$ws.Range("Q1").Value = "Nº Cliente:"
$ws.Range("R1").Value = "Quantidade:"
$ws.Range("S1").Value = "Data Entrega Cliente:"
$ws.Range("T1").Value = "Status, Prazo Compra:"
$ws.Range("U1").Value = "Status Painel:"
$ws.Range("V1").Value = "Reunião com Cliente:"
$ws.Range("W1").Value = "Unnamed: 22"
$ws.Range("X1").Value = "Unnamed: 23"
$ws.Range("Y1").Value = "Data Entrada do Pedido"
$ws.Range("Z1").Value = "Nome"
$ws.Range("AA1").Value = "Idade"
$ws.Range("AB1").Value = "Email"
$ws.Range("AC1").Value = "DataCadastro"
$ws.Range("AD1").Value = "Endereço"
$ws.Range("AE1").Value = "Telefone"
$ws.Range("AF1").Value = "Cidade"
$ws.Range("AG1").Value = "Estado"
$ws.Range("AH1").Value = "País"
$ws.Range("AI1").Value = "CEP"
$ws.Range("AJ1").Value = "Observações"

# --- 2. Existing data rows (2-5): bring the new columns L:AJ into the used
#        range as blank cells (same unstyled format as A2). ----------------
$ws.Range("A2").Copy()
$ws.Range("L2:AJ5").PasteSpecial(-4122)

# --- 3. Brand new rows 6-9: first stamp blank cells across the full A:AJ
#        span so every column is part of the used range, then overwrite the
#        specific cells that carry real values. -----------------------------
$ws.Range("A2").Copy()
$ws.Range("A6:AJ9").PasteSpecial(-4122)

# Row 6
$ws.Range("L6").Value = 10
$ws.Range("M6").Value = 20
$ws.Range("N6").Value = 30
$ws.Range("O6").Value = 40
$ws.Range("P6").Value = "2025-01-29 22:37:43"
$ws.Range("Q6").Value = 50
$ws.Range("R6").Value = 60
$ws.Range("S6").Value = 70
$ws.Range("T6").Value = 80
$ws.Range("U6").Value = 80
$ws.Range("V6").Value = 1001

# Row 7
$ws.Range("L7").Value = 10
$ws.Range("N7").Value = 30
$ws.Range("O7").Value = 40
$ws.Range("P7").Value = "2025-01-29 22:38:58"
$ws.Range("Q7").Value = 50
$ws.Range("R7").Value = 60
$ws.Range("S7").Value = 70
$ws.Range("T7").Value = 80
$ws.Range("U7").Value = 80
$ws.Range("V7").Value = 1001
$ws.Range("W7").Value = 20

# Row 8
$ws.Range("A8").Value = "a"
$ws.Range("J8").Value = "n"
$ws.Range("N8").Value = "x"
$ws.Range("O8").Value = "c"
$ws.Range("Q8").Value = "v"
$ws.Range("R8").Value = "b"
$ws.Range("S8").Value = "n"
$ws.Range("T8").Value = "m"
$ws.Range("U8").Value = ","
$ws.Range("V8").Value = "m"
$ws.Range("X8").Value = "s"
$ws.Range("Y8").Value = "2025-01-29 22:39:56"

# Row 9 (AC9 is a plain date/time text; the rest are digit strings)
$ws.Range("AC9").Value = "2025-01-29 22:41:15"

# --- 4. Digit-like text values: Excel auto-detects plain numeric strings as
#        numbers, so force them to text with a leading apostrophe, then
#        re-paste the (unstyled) A2 format on top so no stray "quote prefix"
#        number format sticks around on the cell's style. -------------------
$textDigitCells = @("J6", "J7", "Z9", "AA9", "AB9", "AD9", "AE9", "AF9", "AG9", "AH9", "AI9", "AJ9")
$textDigitValues = @("90", "90", "1", "2", "3", "4", "5", "6", "7", "8", "9", "0")
for ($i = 0; $i -lt $textDigitCells.Length; $i++) {
    $ws.Range($textDigitCells[$i]).Value = "'" + $textDigitValues[$i]
    $ws.Range("A2").Copy()
    $ws.Range($textDigitCells[$i]).PasteSpecial(-4122)
}
